# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps recorded on the Overview, zh-cn
# and de-de sheets to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 05:04:28"

# --- zh-cn sheet ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 05:04:23"
$wsZhCn.Range("K2").Value = "2016-08-30 05:04:56"

# --- de-de sheet ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 05:04:28"
$wsDeDe.Range("K2").Value = "2016-08-30 05:05:11"
